$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10
$ws.Range("A10").Value = 9751.17
$ws.Range("B10").Value = 9767.7800000000007
$ws.Range("C10").Value = 78.05
$ws.Range("D10").Value = 77.92
$ws.Range("E10").Value = $false
$ws.Range("F10").Value = -0.17
$ws.Range("G10").Value = 42613.765462962961
$ws.Range("H10").Value = $false

# Row 11
$ws.Range("A11").Value = 9831.1299999999992
$ws.Range("B11").Value = 9751.17
$ws.Range("C11").Value = 77.739999999999995
$ws.Range("D11").Value = 78.38
$ws.Range("E11").Value = $false
$ws.Range("F11").Value = 0.82
$ws.Range("G11").Value = 42614.672662037039
$ws.Range("H11").Value = $true
